$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 441
$wsExhibit.Range("F3").Value = 5403
$wsExhibit.Range("F6").Value = 74

# Sheet "全部类型" (all types) - same rows duplicated, update matching cells
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 441
$wsAll.Range("F3").Value = 5403
$wsAll.Range("F7").Value = 74
